$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I20").Value = 'sv'
$ws.Range("J20").Value = 'Statement-opinion'
$ws.Range("I25").Value = 'sv'
$ws.Range("J25").Value = 'Statement-opinion'
$ws.Range("I35").Value = 'sv'
$ws.Range("J35").Value = 'Statement-opinion'
$ws.Range("I39").Value = 'aa'
$ws.Range("J39").Value = 'Agree/Accept'
$ws.Range("I59").Value = 'aa'
$ws.Range("J59").Value = 'Agree/Accept'
$ws.Range("I65").Value = 'aa'
$ws.Range("J65").Value = 'Agree/Accept'
$ws.Range("I84").Value = 'sv'
$ws.Range("J84").Value = 'Statement-opinion'
$ws.Range("I86").Value = 'sv'
$ws.Range("J86").Value = 'Statement-opinion'
$ws.Range("I99").Value = 'sv'
$ws.Range("J99").Value = 'Statement-opinion'
$ws.Range("I102").Value = 'qy'
$ws.Range("J102").Value = 'Yes-No-Question'
$ws.Range("I105").Value = 'aa'
$ws.Range("J105").Value = 'Agree/Accept'
$ws.Range("I106").Value = 'sd'
$ws.Range("J106").Value = 'Statement-non-opinion'
$ws.Range("I108").Value = 'aa'
$ws.Range("J108").Value = 'Agree/Accept'
$ws.Range("I109").Value = 'aa'
$ws.Range("J109").Value = 'Agree/Accept'
$ws.Range("I115").Value = 'ba'
$ws.Range("J115").Value = 'Appreciation'
$ws.Range("I125").Value = 'sv'
$ws.Range("J125").Value = 'Statement-opinion'
$ws.Range("I150").Value = 'ba'
$ws.Range("J150").Value = 'Appreciation'
$ws.Range("I163").Value = 'sv'
$ws.Range("J163").Value = 'Statement-opinion'
$ws.Range("I166").Value = 'b'
$ws.Range("J166").Value = 'Acknowledge (Backchannel)'
$ws.Range("I189").Value = 'sd'
$ws.Range("J189").Value = 'Statement-non-opinion'
$ws.Range("I198").Value = 'sv'
$ws.Range("J198").Value = 'Statement-opinion'
$ws.Range("I200").Value = 'sd'
$ws.Range("J200").Value = 'Statement-non-opinion'
$ws.Range("I203").Value = '%'
$ws.Range("J203").Value = 'Uninterpretable'
$ws.Range("I221").Value = 'sd'
$ws.Range("J221").Value = 'Statement-non-opinion'
$ws.Range("I222").Value = 'sv'
$ws.Range("J222").Value = 'Statement-opinion'
$ws.Range("I228").Value = 'sv'
$ws.Range("J228").Value = 'Statement-opinion'
$ws.Range("I229").Value = '%'
$ws.Range("J229").Value = 'Uninterpretable'
$ws.Range("I231").Value = 'sv'
$ws.Range("J231").Value = 'Statement-opinion'
$ws.Range("I248").Value = 'sd'
$ws.Range("J248").Value = 'Statement-non-opinion'
$ws.Range("I263").Value = 'sd'
$ws.Range("J263").Value = 'Statement-non-opinion'
$ws.Range("I264").Value = 'sd'
$ws.Range("J264").Value = 'Statement-non-opinion'
$ws.Range("I267").Value = 'sd'
$ws.Range("J267").Value = 'Statement-non-opinion'
$ws.Range("I269").Value = 'sd'
$ws.Range("J269").Value = 'Statement-non-opinion'
$ws.Range("I294").Value = 'sv'
$ws.Range("J294").Value = 'Statement-opinion'
$ws.Range("I296").Value = 'sv'
$ws.Range("J296").Value = 'Statement-opinion'
$ws.Range("I297").Value = 'sd'
$ws.Range("J297").Value = 'Statement-non-opinion'
$ws.Range("I301").Value = 'sd'
$ws.Range("J301").Value = 'Statement-non-opinion'
$ws.Range("I306").Value = 'ba'
$ws.Range("J306").Value = 'Appreciation'
$ws.Range("I314").Value = '%'
$ws.Range("J314").Value = 'Uninterpretable'
$ws.Range("I319").Value = 'sd'
$ws.Range("J319").Value = 'Statement-non-opinion'
$ws.Range("I326").Value = 'sd'
$ws.Range("J326").Value = 'Statement-non-opinion'
$ws.Range("I346").Value = 'sd'
$ws.Range("J346").Value = 'Statement-non-opinion'
$ws.Range("I373").Value = 'aa'
$ws.Range("J373").Value = 'Agree/Accept'
$ws.Range("I378").Value = 'sd'
$ws.Range("J378").Value = 'Statement-non-opinion'
$ws.Range("I387").Value = 'sd'
$ws.Range("J387").Value = 'Statement-non-opinion'
$ws.Range("I395").Value = 'sv'
$ws.Range("J395").Value = 'Statement-opinion'
$ws.Range("I397").Value = 'aa'
$ws.Range("J397").Value = 'Agree/Accept'
$ws.Range("I398").Value = 'b'
$ws.Range("J398").Value = 'Acknowledge (Backchannel)'
$ws.Range("I405").Value = 'b'
$ws.Range("J405").Value = 'Acknowledge (Backchannel)'
$ws.Range("I411").Value = 'qy'
$ws.Range("J411").Value = 'Yes-No-Question'
